$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1447.8182
$ws.Range("I131").Value = 864.36365
$ws.Range("J131").Value = 2031.2727
$ws.Range("K131").Value = 2593.09095
$ws.Range("L131").Value = 6093.8181
$ws.Range("M131").Value = 2446.90905
$ws.Range("N131").Value = -16173.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4823.6226
$ws.Range("I32").Value = 3685.4773
$ws.Range("K32").Value = 3685.4773
$ws.Range("M32").Value = -3398.4773
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1731.1818
$ws.Range("I61").Value = 1706.125
$ws.Range("K61").Value = 1706.125
$ws.Range("M61").Value = -1494.125
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1141.6
$ws.Range("I74").Value = 577.3570999999999
$ws.Range("K74").Value = 577.3570999999999
$ws.Range("M74").Value = 296.6429000000001
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1141.6
$ws.Range("I77").Value = 577.3570999999999
$ws.Range("K77").Value = 2886.7855
$ws.Range("M77").Value = 1481.2145
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2009.5
$ws.Range("I102").Value = 2009.5
$ws.Range("K102").Value = 2009.5
$ws.Range("M102").Value = -387.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2061.6
$ws.Range("I122").Value = 2048
$ws.Range("J122").Value = 2150
$ws.Range("K122").Value = 6144
$ws.Range("L122").Value = 6450
$ws.Range("M122").Value = -3694
$ws.Range("N122").Value = -11350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1425.5
$ws.Range("I132").Value = 1178.7916
$ws.Range("K132").Value = 3536.3748
$ws.Range("M132").Value = -1006.3748
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1731.1818
$ws.Range("I136").Value = 1706.125
$ws.Range("K136").Value = 5118.375
$ws.Range("M136").Value = -2568.375
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 15333
$ws.Range("I33").Value = 12999.5
$ws.Range("K33").Value = 12999.5
$ws.Range("M33").Value = -12663.5
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3222.8667
$ws.Range("I31").Value = 3997.2
$ws.Range("J31").Value = 2835.7
$ws.Range("K31").Value = 3997.2
$ws.Range("L31").Value = 2835.7
$ws.Range("M31").Value = -3702.2
$ws.Range("N31").Value = -3425.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3222.8667
$ws.Range("I34").Value = 3997.2
$ws.Range("J34").Value = 2835.7
$ws.Range("K34").Value = 3997.2
$ws.Range("L34").Value = 2835.7
$ws.Range("M34").Value = -3795.2
$ws.Range("N34").Value = -3239.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4035.3333
$ws.Range("I58").Value = 2149.6667
$ws.Range("K58").Value = 2149.6667
$ws.Range("M58").Value = -1946.6667
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3732
$ws.Range("I99").Value = 2446
$ws.Range("J99").Value = 4375
$ws.Range("K99").Value = 2446
$ws.Range("L99").Value = 4375
$ws.Range("M99").Value = -948
$ws.Range("N99").Value = -7371

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1090.5
$ws.Range("I105").Value = 1103.4286
$ws.Range("K105").Value = 1103.4286
$ws.Range("M105").Value = 643.5714
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5540.8
$ws.Range("I122").Value = 4056.7144
$ws.Range("K122").Value = 12170.1432
$ws.Range("M122").Value = -9720.143199999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3732
$ws.Range("I126").Value = 2446
$ws.Range("J126").Value = 4375
$ws.Range("K126").Value = 7338
$ws.Range("L126").Value = 13125
$ws.Range("M126").Value = -4868
$ws.Range("N126").Value = -18065

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2667.3462
$ws.Range("I134").Value = 2267.85
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 6803.549999999999
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -4268.549999999999
$ws.Range("N134").Value = -17067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4035.3333
$ws.Range("I136").Value = 2149.6667
$ws.Range("K136").Value = 6449.000100000001
$ws.Range("M136").Value = -3899.000100000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1686.6666
$ws.Range("I3").Value = 1686.6666
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5059.9998
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4947.9998
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 439.85715
$ws.Range("J5").Value = 899
$ws.Range("L5").Value = 2697
$ws.Range("N5").Value = -2921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 26014.75
$ws.Range("I63").Value = 29.5
$ws.Range("J63").Value = 52000
$ws.Range("K63").Value = 88.5
$ws.Range("L63").Value = 156000
$ws.Range("M63").Value = 660.5
$ws.Range("N63").Value = -157498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 26014.75
$ws.Range("I66").Value = 29.5
$ws.Range("J66").Value = 52000
$ws.Range("K66").Value = 265.5
$ws.Range("L66").Value = 468000
$ws.Range("M66").Value = 3478.5
$ws.Range("N66").Value = -475488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2717.6365
$ws.Range("I69").Value = 2249.75
$ws.Range("K69").Value = 6749.25
$ws.Range("M69").Value = -5938.25
$ws.Range("N69").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 2717.6365
$ws.Range("I72").Value = 2249.75
$ws.Range("K72").Value = 20247.75
$ws.Range("M72").Value = -16191.75
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 439.85715
$ws.Range("J135").Value = 899
$ws.Range("L135").Value = 8091
$ws.Range("N135").Value = -13161

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2506.6453
$ws.Range("I140").Value = 1269.6875
$ws.Range("J140").Value = 3826.0667
$ws.Range("K140").Value = 3809.0625
$ws.Range("L140").Value = 11478.2001
$ws.Range("M140").Value = 1370.9375
$ws.Range("N140").Value = -21838.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1200
$ws.Range("J80").Value = 1200
$ws.Range("L80").Value = 1200
$ws.Range("N80").Value = -3196

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1200
$ws.Range("J83").Value = 1200
$ws.Range("L83").Value = 6000
$ws.Range("N83").Value = -15984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4064.6538
$ws.Range("I102").Value = 5307.2144
$ws.Range("J102").Value = 2615
$ws.Range("K102").Value = 5307.2144
$ws.Range("L102").Value = 2615
$ws.Range("M102").Value = -3685.2144
$ws.Range("N102").Value = -5859

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3259.6924
$ws.Range("I132").Value = 2628.5625
$ws.Range("K132").Value = 7885.6875
$ws.Range("M132").Value = -5355.6875
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13856.417
$ws.Range("I40").Value = 15110.375
$ws.Range("K40").Value = 15110.375
$ws.Range("M40").Value = -14974.375
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2053.3333
$ws.Range("I82").Value = 1396.25
$ws.Range("J82").Value = 3367.5
$ws.Range("K82").Value = 1396.25
$ws.Range("L82").Value = 3367.5
$ws.Range("M82").Value = -1035.25
$ws.Range("N82").Value = -4089.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2053.3333
$ws.Range("I85").Value = 1396.25
$ws.Range("J85").Value = 3367.5
$ws.Range("K85").Value = 1396.25
$ws.Range("L85").Value = 3367.5
$ws.Range("M85").Value = -148.25
$ws.Range("N85").Value = -5863.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6657.4287
$ws.Range("I122").Value = 5949.8335
$ws.Range("J122").Value = 7188.125
$ws.Range("K122").Value = 17849.5005
$ws.Range("L122").Value = 21564.375
$ws.Range("M122").Value = -15399.5005
$ws.Range("N122").Value = -26464.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29777.5
$ws.Range("J70").Value = 29777.5
$ws.Range("L70").Value = 29777.5
$ws.Range("N70").Value = -30407.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 29777.5
$ws.Range("J73").Value = 29777.5
$ws.Range("L73").Value = 29777.5
$ws.Range("N73").Value = -31961.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1782.7667
$ws.Range("I132").Value = 1093.9
$ws.Range("J132").Value = 3160.5
$ws.Range("K132").Value = 3281.7
$ws.Range("L132").Value = 9481.5
$ws.Range("M132").Value = -751.7000000000003
$ws.Range("N132").Value = -14541.5
